$wb = $excel.ActiveWorkbook

# Sheet "OFF" - update row 2 (Texans offensive target depth data for Week 13)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 221
$wsOff.Range("C2").Value = 138
$wsOff.Range("D2").Value = 66
$wsOff.Range("F2").Value = 5

# Sheet "DEF" - update row 2 (Texans defensive target depth data for Week 13)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 353
$wsDef.Range("C2").Value = 256
$wsDef.Range("D2").Value = 68
$wsDef.Range("E2").Value = 36
